$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (RandomForestRegressor) - values change, name stays the same
$ws.Range("B3").Value = 0.9977675495182273
$ws.Range("C3").Value = 0.9977804364125892
$ws.Range("D3").Value = 0.9866095193045398

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9969037223432045
$ws.Range("C4").Value = 0.997012005038021
$ws.Range("D4").Value = 0.9680572685078858

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9988769499527667
$ws.Range("C5").Value = 0.9986691538517337
$ws.Range("D5").Value = 0.998339619466513
